$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Open-date values in column D (rows 2-6): 45200 -> 45206 (Oct 1, 2023 -> Oct 7, 2023)
$ws.Range("D2:D6").Value = 45206

# Update the active selection to match the final saved state
$ws.Range("E10").Select()
